$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.593.16'
$ws.Range("E2").Value = '  -1.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.051.46'
$ws.Range("E3").Value = '  -0.91%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.40'
$ws.Range("E5").Value = '  -3.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.71'
$ws.Range("E6").Value = '  -3.09%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.042.90'
$ws.Range("E8").Value = '  -0.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +0.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.12'
$ws.Range("E11").Value = '  -8.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  +0.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  +3.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.92'
$ws.Range("E14").Value = '  -2.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.545.93'
$ws.Range("E15").Value = '  -1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.656.14'
$ws.Range("E16").Value = '  -1.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.110'
$ws.Range("E17").Value = '  -0.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.057.06'
$ws.Range("E18").Value = '  -1.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.59'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.52'
$ws.Range("E20").Value = '  -4.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.16'
$ws.Range("E21").Value = '  -2.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  -1.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.07'
$ws.Range("E23").Value = '  -1.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.84'
$ws.Range("E24").Value = '  +2.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.00'
$ws.Range("E25").Value = '  -1.98%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("E27").Value = '  -2.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.02'
$ws.Range("E28").Value = '  -2.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.81'
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.84'
$ws.Range("E31").Value = '  -8.73%  '

$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '56.47'
$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.32'
$ws.Range("E34").Value = '  -7.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.29'
$ws.Range("E35").Value = '  +3.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.91'
$ws.Range("E36").Value = '  +1.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '471.93'
$ws.Range("E37").Value = '  -10.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0392'
$ws.Range("E38").Value = '  -4.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.077.31'
$ws.Range("E39").Value = '  +0.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0789'
$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.114'
$ws.Range("E41").Value = '  -1.90%  '

$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.02'
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.62'
$ws.Range("E43").Value = '  +2.41%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.249'
$ws.Range("E45").Value = '  -0.50%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₃0535'
$ws.Range("E46").Value = '  +9.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.92'
$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("E48").Value = '  -3.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.28'
$ws.Range("E49").Value = '  +2.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.107'
$ws.Range("E50").Value = '  +1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.28'
$ws.Range("E51").Value = '  +1.90%  '
